# Changes made to plot.
# Rescale the raw CFU counts (column B, rows 2-9) on every sheet by dividing
# by 20, then update the saved selection/active-sheet view state to match.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1", "2", "3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($row = 2; $row -le 9; $row++) {
        $cell = $ws.Cells.Item($row, 2)
        $cell.Value2 = $cell.Value2 / 20
    }
}

# Update each sheet's stored selection.
$ws1 = $wb.Worksheets.Item("1")
$ws2 = $wb.Worksheets.Item("2")
$ws3 = $wb.Worksheets.Item("3")

[void]$ws1.Range("E3").Select()
[void]$ws2.Range("D3").Select()

# Sheet "3" ends up active/selected with cell D1 selected.
[void]$ws3.Activate()
[void]$ws3.Range("D1").Select()
